$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.360.33'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.687.62'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '680.13'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.23'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.05'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.28%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  -3.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.311.87'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.45'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.683.53'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.333.83'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.89'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.00'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.86'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.833.40'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -5.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.97'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.15'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.83%  '
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.63'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.01%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.00'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.678.23'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.25'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.27%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -2.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '170.31'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.93%  '
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.70'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.37%  '
$ws.Range('E46').Value = '  -1.98%  '
$ws.Range('E47').Value = '  -4.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.72'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.95%  '
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('E50').Value = '  -4.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.80'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.70%  '
